$wb = $excel.ActiveWorkbook

# Row number (F column) -> new value, to be applied on both the
# "展览" and "全部类型" worksheets.
$updates = @{
    2  = 1543
    3  = 43
    4  = 1004
    5  = 10
    7  = 2521
    9  = 1567
    11 = 179
    12 = 61
    13 = 477
    15 = 39
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
